$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# Heading label abbreviations: drop the accent on the final vowel.
Replace-All "A = ASEPTÁ" "A = ASEPTA"
Replace-All "E = EKSPLORÁ" "E = EKSPLORA"
Replace-All "K = KONEKTÁ" "K = KONEKTA"
Replace-All "P = PRAKTIKÁ" "P = PRAKTIKA"

# Paragraph after the A-E-K-P bullet list.
Replace-All "Por apliká A-K-E-P na kualke momentu durante e seshon, praktikando abilidatnan i durante reflekshon." "Por apliká A-K-E-P na kualke momentu durante e seshon, segun bo ta praktiká abilidatnan i durante reflekshon."

# "Gradisí..." bullet.
Replace-All "Gradisí i rekonosé nan kontribushonnan, spesialmente na momentu ku kompartí ta difísil." "Gradisí i rekonosé nan kontribushonnan, spesialmente ora ku ta difísil pa kompartí ."

# "Ehèmpel for di un Práktika di Grupo..." line.
Replace-All "Ehèmpel for di un Práktika di Grupo (praktikando Tempu abo ku bo yu):  " "Ehèmpel di un Práktika di Grupo (praktiká Tempu pa abo ku bo yu):  "

# "Usa afirmashonnan manera..." bullet.
Replace-All "Usa afirmashonnan manera: “Esei ta un tremendo ehèmpel di krea konfiansa ku bo yu.”" "Usa afirmashonnan manera: “Esei ta un tremendo ehèmpel kon pa krea konfiansa ku bo yu.”"

# "Resumen..." paragraph: "tokante Tempu abo ku bo Yu" -> "tokante Tempu pa abo ku bo Yu"
Replace-All "tokante Tempu abo ku bo Yu.) " "tokante Tempu pa abo ku bo Yu.) "

# "Esaki ta nifiká..." paragraph in the P = PRAKTIKA row.
Replace-All "Esaki ta nifiká ku bo mester duna mayornan mas tantu oportunidat posibel pa praktiká abilidatnan klave, sea den Práktika den Grupo òf Práktika den Par." "Esaki ta nifiká ku bo mester duna mayornan mas tantu oportunidat posibel pa praktiká abilidatnan klave, sea den Grupo òf den Par."
